$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update match rows whose betexplorer ordering changed (2023-11-24, 2023-12-09, 2023-12-31 matchdays) ---

# Row 74
$ws.Range("F74").Value = "Aluminium Arak"
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = "Paykan"
$ws.Range("I74").Value = 1
$ws.Range("J74").Value = 1.85
$ws.Range("K74").Value = "22/11/2023 15:12"
$ws.Range("L74").Value = 2.1
$ws.Range("M74").Value = "24/11/2023 12:10"
$ws.Range("N74").Value = 2.8
$ws.Range("O74").Value = "22/11/2023 15:12"
$ws.Range("P74").Value = 2.36
$ws.Range("Q74").Value = "24/11/2023 12:10"
$ws.Range("R74").Value = 4.58
$ws.Range("S74").Value = "22/11/2023 15:12"
$ws.Range("T74").Value = 4.43
$ws.Range("U74").Value = "24/11/2023 12:10"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/aluminium-arak-paykan/2o67LExL/"

# Row 75
$ws.Range("F75").Value = "Havadar SC"
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = "Malavan"
$ws.Range("I75").Value = 1
$ws.Range("J75").Value = 2.51
$ws.Range("K75").Value = "22/11/2023 15:12"
$ws.Range("L75").Value = 3.16
$ws.Range("M75").Value = "24/11/2023 11:20"
$ws.Range("N75").Value = 2.65
$ws.Range("O75").Value = "22/11/2023 15:12"
$ws.Range("P75").Value = 2.5
$ws.Range("Q75").Value = "24/11/2023 11:20"
$ws.Range("R75").Value = 2.96
$ws.Range("S75").Value = "22/11/2023 15:12"
$ws.Range("T75").Value = 2.81
$ws.Range("U75").Value = "24/11/2023 11:33"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/havadar-sc-malavan/tx5BKYiR/"

# Row 76
$ws.Range("F76").Value = "Mes Rafsanjan"
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = "Sanat Naft"
$ws.Range("I76").Value = 2
$ws.Range("J76").Value = 1.53
$ws.Range("K76").Value = "22/11/2023 15:12"
$ws.Range("L76").Value = 1.58
$ws.Range("M76").Value = "24/11/2023 12:25"
$ws.Range("N76").Value = 3.43
$ws.Range("O76").Value = "22/11/2023 15:12"
$ws.Range("P76").Value = 3.36
$ws.Range("Q76").Value = "24/11/2023 12:25"
$ws.Range("R76").Value = 6.31
$ws.Range("S76").Value = "22/11/2023 15:12"
$ws.Range("T76").Value = 7.1
$ws.Range("U76").Value = "24/11/2023 12:25"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/mes-rafsanjan-sanat-naft/MVnvGh0r/"

# Row 84
$ws.Range("F84").Value = "Nassaji Mazandaran"
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = "Esteghlal Khuzestan"
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1.84
$ws.Range("K84").Value = "08/12/2023 00:43"
$ws.Range("L84").Value = 2
$ws.Range("M84").Value = "09/12/2023 12:25"
$ws.Range("N84").Value = 2.89
$ws.Range("O84").Value = "08/12/2023 00:43"
$ws.Range("P84").Value = 2.68
$ws.Range("Q84").Value = "09/12/2023 12:25"
$ws.Range("R84").Value = 4.6
$ws.Range("S84").Value = "08/12/2023 00:43"
$ws.Range("T84").Value = 5.06
$ws.Range("U84").Value = "09/12/2023 11:53"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/mazandaran-esteghlal-khuzestan/tfMGTBUs/"

# Row 85
$ws.Range("F85").Value = "Sepahan"
$ws.Range("G85").Value = 4
$ws.Range("H85").Value = "Mes Rafsanjan"
$ws.Range("I85").Value = 1
$ws.Range("J85").Value = 1.54
$ws.Range("K85").Value = "08/12/2023 00:43"
$ws.Range("L85").Value = 1.6
$ws.Range("M85").Value = "09/12/2023 12:17"
$ws.Range("N85").Value = 3.43
$ws.Range("O85").Value = "08/12/2023 00:43"
$ws.Range("P85").Value = 3.46
$ws.Range("Q85").Value = "09/12/2023 12:17"
$ws.Range("R85").Value = 5.99
$ws.Range("S85").Value = "08/12/2023 00:43"
$ws.Range("T85").Value = 6.48
$ws.Range("U85").Value = "09/12/2023 12:17"
$ws.Range("V85").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/sepahan-mes-rafsanjan/zHLKSVql/"

# Row 112
$ws.Range("F112").Value = "Aluminium Arak"
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = "Malavan"
$ws.Range("I112").Value = 1
$ws.Range("J112").Value = 2.52
$ws.Range("K112").Value = "30/12/2023 00:42"
$ws.Range("L112").Value = 2.87
$ws.Range("M112").Value = "31/12/2023 12:13"
$ws.Range("N112").Value = 2.58
$ws.Range("O112").Value = "30/12/2023 00:42"
$ws.Range("P112").Value = 2.34
$ws.Range("Q112").Value = "31/12/2023 12:13"
$ws.Range("R112").Value = 3.11
$ws.Range("S112").Value = "30/12/2023 00:42"
$ws.Range("T112").Value = 3.39
$ws.Range("U112").Value = "31/12/2023 12:18"
$ws.Range("V112").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/aluminium-arak-malavan/xQBvyFg5/"

# Row 113
$ws.Range("F113").Value = "Havadar SC"
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = "Esteghlal Khuzestan"
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2.24
$ws.Range("K113").Value = "30/12/2023 00:42"
$ws.Range("L113").Value = 2.32
$ws.Range("M113").Value = "31/12/2023 12:26"
$ws.Range("N113").Value = 2.69
$ws.Range("O113").Value = "30/12/2023 00:42"
$ws.Range("P113").Value = 2.61
$ws.Range("Q113").Value = "31/12/2023 12:26"
$ws.Range("R113").Value = 3.48
$ws.Range("S113").Value = "30/12/2023 00:42"
$ws.Range("T113").Value = 3.9
$ws.Range("U113").Value = "31/12/2023 12:26"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/havadar-sc-esteghlal-khuzestan/bZCzxeva/"

# --- Append new match row 116 (Zob Ahan vs Gol Gohar) ---
$ws.Range("A115:V115").Copy()
$ws.Range("A116:V116").PasteSpecial(-4122)
$ws.Range("A116").Value = 115
$ws.Range("B116").Value = "iran"
$ws.Range("C116").Value = "persian-gulf-pro-league"
$ws.Range("D116").Value = "2023-2024"
$ws.Range("E116").Value = 45294.52083333334
$ws.Range("F116").Value = "Zob Ahan"
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = "Gol Gohar"
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 2.39
$ws.Range("K116").Value = "02/01/2024 09:12"
$ws.Range("L116").Value = 2.52
$ws.Range("M116").Value = "03/01/2024 12:29"
$ws.Range("N116").Value = 2.58
$ws.Range("O116").Value = "02/01/2024 09:12"
$ws.Range("P116").Value = 2.43
$ws.Range("Q116").Value = "03/01/2024 12:29"
$ws.Range("R116").Value = 3.36
$ws.Range("S116").Value = "02/01/2024 09:12"
$ws.Range("T116").Value = 3.79
$ws.Range("U116").Value = "03/01/2024 12:29"
$ws.Range("V116").Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/zob-ahan-gol-gohar/IZJnZhOH/"
